# Correção nos dados: a linha 6 ("grandes regiões e unidades da federação")
# era apenas um rótulo de seção sem valores - os dados da região "norte"
# estavam, por engano, uma linha abaixo do seu rótulo. Remove essa linha de
# rótulo inteira, o que desloca os dados de B:G das linhas 7..37 para
# 6..36 (mantendo os rótulos corretos na coluna A) e elimina a última linha
# (antiga 37) que sobrava.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Delete()
